$d = $word.ActiveDocument

# Locate the target paragraph (the sextortion paragraph) by its distinctive leading text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*PsiXBot malware has been updated*") {
        $target = $p
        break
    }
}

# Helper that replaces the text found between two anchor substrings (relative
# to the *current* paragraph text) with new text, forcing the replacement to
# live in its own run (so the resulting OOXML mirrors a manual, word-by-word
# revision instead of one big run rewrite). Edits are applied strictly
# left-to-right across the paragraph so earlier, already-created run
# boundaries are never disturbed by a later edit.
function Set-RunText($searchAnchor, $offsetInAnchor, $oldLength, $newText) {
    $pStart = $target.Range.Start
    $pText = $target.Range.Text
    $idx = $pText.IndexOf($searchAnchor)
    if ($idx -lt 0) {
        throw "anchor not found: $searchAnchor"
    }
    $relStart = $idx + $offsetInAnchor
    $relEnd = $relStart + $oldLength
    $r = $d.Range($pStart + $relStart, $pStart + $relEnd)
    $r.Font.Bold = $true
    $r.Text = $newText
    $r.Font.Bold = $false
}

# 1) "you" -> "anyone"   (... true. If you have visited adult websites ...)
Set-RunText "If you have visited" 3 3 "anyone"

# 2) "ve" -> "s"          (have -> has)
Set-RunText "have visited" 2 2 "s"

# 3) "you" -> "they"     (credibility that you have been recorded)
Set-RunText "that you have been recorded" 5 3 "they"

# 4) "you've" -> "anyone" (However, if you've received)
Set-RunText "if you've received" 3 6 "anyone"

# 5) "you" -> "them"      (blackmail you.)
Set-RunText "blackmail you." 10 3 "them"

# 6) "you" -> "they"      (even if you agree to pay)
Set-RunText "even if you agree" 8 3 "they"

# 7) "you" -> "nobody"    (For this reason, you should)
Set-RunText "reason, you should" 8 3 "nobody"

# 8) "never " -> "" (deleted; "you should never trust" -> "nobody should trust")
Set-RunText "nobody should never trust" 14 6 ""

Write-Output "done"
